# Update countries & provincias Spain
# - Swap the row152/row153 country labels: Georgia's case counts overtook
#   "Republica de Chipre", so Georgia now appears in row 152 (with fresh,
#   larger numbers) and "Republica de Chipre" drops to row 153 (keeping the
#   numbers that used to belong to row 152).
# - Refresh the daily COVID figures (Casos totales, Nuevos casos, Casos
#   activos, Recuperados, Muertes hoy, Muertes) for several countries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# India (row 6)
$ws.Range("B6").Value = 3694878
$ws.Range("C6").Value = 6939
$ws.Range("D6").Value = 2840040
$ws.Range("E6").Value = 789369
$ws.Range("G6").Value = 34
$ws.Range("H6").Value = 65469

# Kazajistan (row 33)
$ws.Range("D33").Value = 96899
$ws.Range("E33").Value = 7450

# Armenia (row 60)
$ws.Range("B60").Value = 43878
$ws.Range("C60").Value = 97
$ws.Range("D60").Value = 38356
$ws.Range("E60").Value = 4641
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 881

# Uzbekistan (row 62)
$ws.Range("B62").Value = 41994
$ws.Range("C62").Value = 101
$ws.Range("E62").Value = 2397
$ws.Range("G62").Value = 2
$ws.Range("H62").Value = 322

# Australia (row 72)
$ws.Range("D72").Value = 21503
$ws.Range("E72").Value = 3659

# Hungria (row 107)
$ws.Range("B107").Value = 6257
$ws.Range("C107").Value = 118
$ws.Range("D107").Value = 3821
$ws.Range("E107").Value = 1820
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 616

# Row 152 now becomes Georgia with its updated, higher figures
$ws.Range("A152").Value = "Georgia"
$ws.Range("B152").Value = 1510
$ws.Range("C152").Value = 23
$ws.Range("D152").Value = 1243
$ws.Range("E152").Value = 248
$ws.Range("H152").Value = 19

# Row 153 now becomes Republica de Chipre, carrying the figures that used
# to sit in row 152 before Georgia's update overtook it
$ws.Range("A153").Value = "Republica de Chipre"
$ws.Range("B153").Value = 1488
$ws.Range("D153").Value = 1139
$ws.Range("E153").Value = 329
$ws.Range("H153").Value = 20
